# Update cryptos list figures (price + volume change) per data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.076.49"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "3.752.81"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'604.28"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'169.35"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").Value = "3.751.42"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("E10").Value = "  +5.68%  "
$ws.Range("D11").Value = "'6.37"
$ws.Range("E11").Value = "  +2.83%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'38.38"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").Value = "'0.0000250"
$ws.Range("E14").Value = "  +4.39%  "
$ws.Range("D15").Value = "4.379.03"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "3.756.24"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "69.076.44"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D20").Value = "'17.11"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "'10.84"
$ws.Range("E21").Value = "  +18.66%  "
$ws.Range("D22").Value = "'494.64"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'0.729"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").Value = "'0.0000156"
$ws.Range("E24").Value = "  +12.98%  "
$ws.Range("D25").Value = "'85.50"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("D27").Value = "'12.37"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").Value = "'10.32"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +7.67%  "
$ws.Range("D31").Value = "'2.99"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "'7.94"
$ws.Range("E32").Value = "  +4.12%  "
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").Value = "3.899.21"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "3.686.97"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").Value = "'5.89"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").Value = "'0.323"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").Value = "'3.08"
$ws.Range("E42").Value = "  +10.69%  "
$ws.Range("D43").Value = "'437.81"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "'48.56"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").Value = "'8.47"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D48").Value = "'40.40"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "'141.09"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.801.97"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0356"
$ws.Range("E51").Value = "  +2.55%  "
